$wb = $excel.ActiveWorkbook

$oldGuid = "d1694ec7-617f-4450-a1ce-9ec8fe74f91f"
$newGuid = "e8789362-58e0-4717-9011-f7d29894228c"
$oldHash = "d9bc6ad6e8e8e7f661545180268eb96cfe818462"
$newHash = "bd19df8fab97ac6815e8d0b6c74b1e36537adddd"

$oldMdName   = "$oldGuid.md"
$newMdName   = "$newGuid.md"
$oldZhXlf    = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf    = "$newGuid.$newHash.zh-cn.xlf"
$oldDeXlf    = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf    = "$newGuid.$newHash.de-de.xlf"

$oldZhTime = "2016-03-08 05:26:49"
$newZhTime = "2016-03-08 05:27:33"
$oldDeTime = "2016-03-08 05:26:59"
$newDeTime = "2016-03-08 05:27:43"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/5a203a1741d5f99c0494ec1b546b6f6e6795e35c/e2e/$oldMdName"
$cfgUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/5a203a1741d5f99c0494ec1b546b6f6e6795e35c/.localization-config"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/678402096051ac330d9d697c2ae44c06a57883ff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$oldZhXlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/439820ba6974af8e27e0e33ccfbf9cecb1386382/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$oldDeXlf"

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhTime
$wsZh.Cells.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", $newZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeTime
$wsDe.Cells.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", $newDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

"done"
